$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 338.6742503333333
$ws.Range("H2").Value = 1016.022751
$ws.Range("I2").Value = 0.5849329800180821
$ws.Range("J2").Value = 0.584932980018082
$ws.Range("M2").Value = 145.7007446666667
$ws.Range("N2").Value = 437.1022340000001
$ws.Range("O2").Value = 0.2865937750105843
$ws.Range("P2").Value = 0.2865937750105843
$ws.Range("Q2").Value = 49345.09047299175
$ws.Range("R2").Value = 444105.8142569258
$ws.Range("S2").Value = 0.1676381508715728
$ws.Range("T2").Value = 0.1676381508715728
$ws.Range("G3").Value = 338.6742503333333
$ws.Range("H3").Value = 1016.022751
$ws.Range("I3").Value = 0.5849329800180821
$ws.Range("J3").Value = 0.584932980018082
$ws.Range("O3").Value = 0.3320294904365841
$ws.Range("P3").Value = 0.3320294904365841
$ws.Range("Q3").Value = 57168.1127571229
$ws.Range("R3").Value = 514513.0148141062
$ws.Range("S3").Value = 0.1942149992949565
$ws.Range("T3").Value = 0.1942149992949564
$ws.Range("G4").Value = 338.6742503333333
$ws.Range("H4").Value = 1016.022751
$ws.Range("I4").Value = 0.5849329800180821
$ws.Range("J4").Value = 0.584932980018082
$ws.Range("M4").Value = 128.1261546666667
$ws.Range("N4").Value = 384.378464
$ws.Range("O4").Value = 0.2520245069956105
$ws.Range("P4").Value = 0.2520245069956105
$ws.Range("Q4").Value = 43393.02937982605
$ws.Range("R4").Value = 390537.2644184345
$ws.Range("S4").Value = 0.1474174459145304
$ws.Range("T4").Value = 0.1474174459145304
$ws.Range("G5").Value = 338.6742503333333
$ws.Range("H5").Value = 1016.022751
$ws.Range("I5").Value = 0.5849329800180821
$ws.Range("J5").Value = 0.584932980018082
$ws.Range("M5").Value = 65.761079
$ws.Range("N5").Value = 197.283237
$ws.Range("O5").Value = 0.1293522275572212
$ws.Range("P5").Value = 0.1293522275572212
$ws.Range("Q5").Value = 22271.58413143611
$ws.Range("R5").Value = 200444.257182925
$ws.Range("S5").Value = 0.0756623839370225
$ws.Range("T5").Value = 0.07566238393702249
$ws.Range("I6").Value = 0.279688040971731
$ws.Range("J6").Value = 0.2796880409717309
$ws.Range("M6").Value = 145.7007446666667
$ws.Range("N6").Value = 437.1022340000001
$ws.Range("O6").Value = 0.2865937750105843
$ws.Range("P6").Value = 0.2865937750105843
$ws.Range("Q6").Value = 23594.55212379587
$ws.Range("R6").Value = 212350.9691141629
$ws.Range("S6").Value = 0.08015685148740334
$ws.Range("T6").Value = 0.08015685148740333
$ws.Range("I7").Value = 0.279688040971731
$ws.Range("J7").Value = 0.2796880409717309
$ws.Range("O7").Value = 0.3320294904365841
$ws.Range("P7").Value = 0.3320294904365841
$ws.Range("S7").Value = 0.09286467772505032
$ws.Range("T7").Value = 0.09286467772505029
$ws.Range("I8").Value = 0.279688040971731
$ws.Range("J8").Value = 0.2796880409717309
$ws.Range("M8").Value = 128.1261546666667
$ws.Range("N8").Value = 384.378464
$ws.Range("O8").Value = 0.2520245069956105
$ws.Range("P8").Value = 0.2520245069956105
$ws.Range("Q8").Value = 20748.55033596693
$ws.Range("R8").Value = 186736.9530237024
$ws.Range("S8").Value = 0.0704882406384686
$ws.Range("T8").Value = 0.07048824063846859
$ws.Range("I9").Value = 0.279688040971731
$ws.Range("J9").Value = 0.2796880409717309
$ws.Range("M9").Value = 65.761079
$ws.Range("N9").Value = 197.283237
$ws.Range("O9").Value = 0.1293522275572212
$ws.Range("P9").Value = 0.1293522275572212
$ws.Range("Q9").Value = 10649.24691862287
$ws.Range("R9").Value = 95843.22226760586
$ws.Range("S9").Value = 0.03617827112080876
$ws.Range("T9").Value = 0.03617827112080876
$ws.Range("G10").Value = 77.79536166666666
$ws.Range("H10").Value = 233.386085
$ws.Range("I10").Value = 0.1343623634996766
$ws.Range("J10").Value = 0.1343623634996766
$ws.Range("M10").Value = 145.7007446666667
$ws.Range("N10").Value = 437.1022340000001
$ws.Range("O10").Value = 0.2865937750105843
$ws.Range("P10").Value = 0.2865937750105843
$ws.Range("Q10").Value = 11334.84212644599
$ws.Range("R10").Value = 102013.5791380139
$ws.Range("S10").Value = 0.03850741697471666
$ws.Range("T10").Value = 0.03850741697471665
$ws.Range("G11").Value = 77.79536166666666
$ws.Range("H11").Value = 233.386085
$ws.Range("I11").Value = 0.1343623634996766
$ws.Range("J11").Value = 0.1343623634996766
$ws.Range("O11").Value = 0.3320294904365841
$ws.Range("P11").Value = 0.3320294904365841
$ws.Range("Q11").Value = 13131.83391817913
$ws.Range("R11").Value = 118186.5052636122
$ws.Range("S11").Value = 0.04461226708665272
$ws.Range("T11").Value = 0.04461226708665271
$ws.Range("G12").Value = 77.79536166666666
$ws.Range("H12").Value = 233.386085
$ws.Range("I12").Value = 0.1343623634996766
$ws.Range("J12").Value = 0.1343623634996766
$ws.Range("M12").Value = 128.1261546666667
$ws.Range("N12").Value = 384.378464
$ws.Range("O12").Value = 0.2520245069956105
$ws.Range("P12").Value = 0.2520245069956105
$ws.Range("Q12").Value = 9967.620541252605
$ws.Range("R12").Value = 89708.58487127344
$ws.Range("S12").Value = 0.03386260841977101
$ws.Range("T12").Value = 0.033862608419771
$ws.Range("G13").Value = 77.79536166666666
$ws.Range("H13").Value = 233.386085
$ws.Range("I13").Value = 0.1343623634996766
$ws.Range("J13").Value = 0.1343623634996766
$ws.Range("M13").Value = 65.761079
$ws.Range("N13").Value = 197.283237
$ws.Range("O13").Value = 0.1293522275572212
$ws.Range("P13").Value = 0.1293522275572212
$ws.Range("Q13").Value = 5115.906924395238
$ws.Range("R13").Value = 46043.16231955714
$ws.Range("S13").Value = 0.01738007101853625
$ws.Range("T13").Value = 0.01738007101853624
$ws.Range("G14").Value = 0.5886170000000001
$ws.Range("H14").Value = 1.765851
$ws.Range("I14").Value = 0.001016615510510267
$ws.Range("J14").Value = 0.001016615510510266
$ws.Range("M14").Value = 145.7007446666667
$ws.Range("N14").Value = 437.1022340000001
$ws.Range("O14").Value = 0.2865937750105843
$ws.Range("P14").Value = 0.2865937750105843
$ws.Range("Q14").Value = 85.76193522345935
$ws.Range("R14").Value = 771.8574170111341
$ws.Range("S14").Value = 0.0002913556768914496
$ws.Range("T14").Value = 0.0002913556768914496
$ws.Range("G15").Value = 0.5886170000000001
$ws.Range("H15").Value = 1.765851
$ws.Range("I15").Value = 0.001016615510510267
$ws.Range("J15").Value = 0.001016615510510266
$ws.Range("O15").Value = 0.3320294904365841
$ws.Range("P15").Value = 0.3320294904365841
$ws.Range("Q15").Value = 99.35837458454536
$ws.Range("R15").Value = 894.2253712609081
$ws.Range("S15").Value = 0.0003375463299246517
$ws.Range("T15").Value = 0.0003375463299246516
$ws.Range("G16").Value = 0.5886170000000001
$ws.Range("H16").Value = 1.765851
$ws.Range("I16").Value = 0.001016615510510267
$ws.Range("J16").Value = 0.001016615510510266
$ws.Range("M16").Value = 128.1261546666667
$ws.Range("N16").Value = 384.378464
$ws.Range("O16").Value = 0.2520245069956105
$ws.Range("P16").Value = 0.2520245069956105
$ws.Range("Q16").Value = 75.41723278142935
$ws.Range("R16").Value = 678.755095032864
$ws.Range("S16").Value = 0.0002562120228404408
$ws.Range("T16").Value = 0.0002562120228404408
$ws.Range("G17").Value = 0.5886170000000001
$ws.Range("H17").Value = 1.765851
$ws.Range("I17").Value = 0.001016615510510267
$ws.Range("J17").Value = 0.001016615510510266
$ws.Range("M17").Value = 65.761079
$ws.Range("N17").Value = 197.283237
$ws.Range("O17").Value = 0.1293522275572212
$ws.Range("P17").Value = 0.1293522275572212
$ws.Range("Q17").Value = 38.708089037743
$ws.Range("R17").Value = 348.372801339687
$ws.Range("S17").Value = 0.0001315014808537246
$ws.Range("T17").Value = 0.0001315014808537246

Write-Output "Updated 174 cells with new TPM-derived values"
